$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.83"
$ws.Range("D3").Value = "'32.73"
$ws.Range("E3").Value = "'4.02%"
$ws.Range("D4").Value = "'4.956"
$ws.Range("E4").Value = "'-3.09%"
$ws.Range("E5").Value = "'-1.10%"
$ws.Range("D6").Value = "'1.938"
$ws.Range("E6").Value = "'-17.52%"
$ws.Range("D7").Value = "'7.870"
$ws.Range("E7").Value = "'0.85%"
$ws.Range("D8").Value = "'3.795"
$ws.Range("E8").Value = "'-1.78%"
$ws.Range("D9").Value = "'0.9273"
$ws.Range("E9").Value = "'0.42%"
$ws.Range("D10").Value = "'0.1771"
$ws.Range("E10").Value = "'1.11%"
$ws.Range("D11").Value = "'0.07783"
$ws.Range("E11").Value = "'2.82%"
$ws.Range("D12").Value = "'0.08703"
$ws.Range("E12").Value = "'-6.71%"
$ws.Range("D13").Value = "'0.03147"
$ws.Range("E13").Value = "'4.54%"
$ws.Range("D14").Value = "'0.1003"
$ws.Range("E14").Value = "'0.05%"
$ws.Range("D15").Value = "'0.001523"
$ws.Range("E15").Value = "'0.91%"
$ws.Range("D16").Value = "'0.005776"
$ws.Range("E16").Value = "'-2.90%"
$ws.Range("D17").Value = "'3.461"
$ws.Range("E17").Value = "'-0.39%"
$ws.Range("D18").Value = "'2.154"
$ws.Range("E18").Value = "'-4.88%"
$ws.Range("D19").Value = "'0.3312"
$ws.Range("E19").Value = "'1.24%"
$ws.Range("E20").Value = "'0.73%"
$ws.Range("D21").Value = "'4.293"
$ws.Range("E21").Value = "'10.15%"
$ws.Range("E22").Value = "'17.06%"
$ws.Range("D23").Value = "'0.04572"
$ws.Range("E23").Value = "'-1.14%"
$ws.Range("E24").Value = "'-2.14%"
$ws.Range("D25").Value = "'0.004441"
$ws.Range("E25").Value = "'-0.66%"
$ws.Range("E26").Value = "'4.21%"
$ws.Range("D39").Value = "'0.01712"
$ws.Range("E39").Value = "'-1.85%"
$ws.Range("D40").Value = "'0.04765"
$ws.Range("E40").Value = "'3.32%"
$ws.Range("D41").Value = "'0.007511"
$ws.Range("E41").Value = "'7.82%"
$ws.Range("D42").Value = "'0.1355"
$ws.Range("E42").Value = "'-0.49%"
$ws.Range("E43").Value = "'5.45%"
$ws.Range("D44").Value = "'0.01164"
$ws.Range("E44").Value = "'12.87%"
$ws.Range("D45").Value = "'0.00006260"
$ws.Range("E45").Value = "'-0.39%"
$ws.Range("E46").Value = "'0.07%"
$ws.Range("E47").Value = "'-61.13%"
$ws.Range("D48").Value = "'0.8204"
$ws.Range("E48").Value = "'-29.01%"
$ws.Range("E49").Value = "'0.07%"
$ws.Range("E50").Value = "'0.07%"
